$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Constant attributes shared by every "Comercializadora del Agro de Limari /
# Poroto granado" record in this block (columns A,B,C,E,F,G,H,I,O,R never
# change row to row here).
$marketId = 2
$market   = "Comercializadora del Agro de Limarí"
$region   = "Coquimbo"
$codreg   = 4
$catId    = 100112030
$cat      = "Poroto granado"
$variedad = "Sin especificar"
$calidad  = "Primera"
$origen   = "Provincia de Limarí"
$clasif   = "Hortaliza"

function Set-Registro($Row, $Fecha, $Volumen, $PrecioMin, $PrecioMax, $PrecioProm, $Unidad, $PrecioKg, $KgUnidades) {
    $ws.Cells.Item($Row, 1).Value  = $marketId
    $ws.Cells.Item($Row, 2).Value  = $market
    $ws.Cells.Item($Row, 3).Value  = $region
    $ws.Cells.Item($Row, 4).Value  = $Fecha
    $ws.Cells.Item($Row, 5).Value  = $codreg
    $ws.Cells.Item($Row, 6).Value  = $catId
    $ws.Cells.Item($Row, 7).Value  = $cat
    $ws.Cells.Item($Row, 8).Value  = $variedad
    $ws.Cells.Item($Row, 9).Value  = $calidad
    $ws.Cells.Item($Row, 10).Value = $Volumen
    $ws.Cells.Item($Row, 11).Value = $PrecioMin
    $ws.Cells.Item($Row, 12).Value = $PrecioMax
    $ws.Cells.Item($Row, 13).Value = $PrecioProm
    $ws.Cells.Item($Row, 14).Value = $Unidad
    $ws.Cells.Item($Row, 15).Value = $origen
    $ws.Cells.Item($Row, 16).Value = $PrecioKg
    $ws.Cells.Item($Row, 17).Value = $KgUnidades
    $ws.Cells.Item($Row, 18).Value = $clasif
}

# A new weekly record is inserted at the top of the date-sorted block (row
# 83), pushing every existing record down by one row.
$ws.Rows.Item(83).Insert()
Set-Registro 83 45007 600 23000 25000 24000 "$/malla 25 kilos" 960 25

# A second new weekly record is inserted further down the block (at what is
# now row 108), pushing the remaining older records down by one more row.
$ws.Rows.Item(108).Insert()
Set-Registro 108 45008 600 24000 25000 24500 "$/malla 25 kilos" 980 25
